# Auto-generated edit script for 'Add data for 2022-12-04' commit
# Applies +1 corrections to YTD crime-count cells across several sheets
# (Citywide Totals, By Neighborhood, and 9 individual neighborhood sheets),
# reflecting a backfilled/late-reported incident.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("E3").Value = 144
$ws.Range("B6").Value = 374
$ws.Range("C6").Value = 475
$ws.Range("D6").Value = 413
$ws.Range("E6").Value = 469
$ws.Range("F6").Value = 531
$ws.Range("I6").Value = 497
$ws.Range("B7").Value = 499
$ws.Range("C7").Value = 630
$ws.Range("D7").Value = 644
$ws.Range("E7").Value = 694
$ws.Range("F7").Value = 767
$ws.Range("I7").Value = 829

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 53
$ws.Range("E7").Value = 66

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 34
$ws.Range("C7").Value = 39

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F7").Value = 11
$ws.Range("B8").Value = 31
$ws.Range("I21").Value = 15
$ws.Range("E32").Value = 66
$ws.Range("C36").Value = 39
$ws.Range("E53").Value = 82
$ws.Range("F53").Value = 82
$ws.Range("D65").Value = 25
$ws.Range("E91").Value = 7
$ws.Range("B96").Value = 16
$ws.Range("B98").Value = 499
$ws.Range("C98").Value = 630
$ws.Range("D98").Value = 644
$ws.Range("E98").Value = 694
$ws.Range("F98").Value = 767
$ws.Range("I98").Value = 829

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 16

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E3").Value = 13
$ws.Range("F6").Value = 61
$ws.Range("E7").Value = 82
$ws.Range("F7").Value = 82

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("E6").Value = 6
$ws.Range("E7").Value = 7

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("D5").Value = 24
$ws.Range("D6").Value = 25

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 11

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("B5").Value = 22
$ws.Range("B6").Value = 31

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 15

Write-Host "Done applying 2022-12-04 data updates."
